$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("H28").Value = 2194.6
$ws.Range("I28").Value = 908.6667
$ws.Range("J28").Value = 7338.3335
$ws.Range("K28").Value = 908.6667
$ws.Range("L28").Value = 7338.3335
$ws.Range("M28").Value = -423.6667
$ws.Range("N28").Value = -8308.333500000001
$ws.Range("H34").Value = 1539.4
$ws.Range("I34").Value = 1539.4
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1539.4
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1336.4
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 1539.4
$ws.Range("I36").Value = 1539.4
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1539.4
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -824.4000000000001
$ws.Range("N36").ClearContents()
$ws.Range("H43").Value = 1229.6666
$ws.Range("I43").Value = 1244.5
$ws.Range("J43").Value = 1200
$ws.Range("K43").Value = 1244.5
$ws.Range("L43").Value = 1200
$ws.Range("M43").Value = -1175.5
$ws.Range("N43").Value = -1338
$ws.Range("H51").Value = 18268.54
$ws.Range("I51").Value = 16317.182
$ws.Range("J51").Value = 29001
$ws.Range("K51").Value = 16317.182
$ws.Range("L51").Value = 29001
$ws.Range("M51").Value = -15833.182
$ws.Range("N51").Value = -29969
$ws.Range("H62").Value = 4748.4
$ws.Range("I62").Value = 5697.8
$ws.Range("J62").Value = 3799
$ws.Range("K62").Value = 5697.8
$ws.Range("L62").Value = 3799
$ws.Range("M62").Value = -5073.8
$ws.Range("N62").Value = -5047
$ws.Range("H65").Value = 4748.4
$ws.Range("I65").Value = 5697.8
$ws.Range("J65").Value = 3799
$ws.Range("K65").Value = 28489
$ws.Range("L65").Value = 18995
$ws.Range("M65").Value = -25369
$ws.Range("N65").Value = -25235
$ws.Range("H70").Value = 2219.389
$ws.Range("I70").Value = 2100.3333
$ws.Range("J70").Value = 2278.9167
$ws.Range("K70").Value = 6300.999899999999
$ws.Range("L70").Value = 6836.750100000001
$ws.Range("M70").Value = -6030.999899999999
$ws.Range("N70").Value = -7376.750100000001
$ws.Range("H73").Value = 2219.389
$ws.Range("I73").Value = 2100.3333
$ws.Range("J73").Value = 2278.9167
$ws.Range("K73").Value = 6300.999899999999
$ws.Range("L73").Value = 6836.750100000001
$ws.Range("M73").Value = -5364.999899999999
$ws.Range("N73").Value = -8708.750100000001
$ws.Range("H74").Value = 18108.75
$ws.Range("I74").Value = 6410
$ws.Range("J74").Value = 100000
$ws.Range("K74").Value = 6410
$ws.Range("L74").Value = 100000
$ws.Range("M74").Value = -5474
$ws.Range("H76").Value = 4477.909
$ws.Range("I76").Value = 4470
$ws.Range("J76").Value = 4499
$ws.Range("K76").Value = 4470
$ws.Range("L76").Value = 4499
$ws.Range("M76").Value = -4155
$ws.Range("N76").Value = -5129
$ws.Range("H77").Value = 18108.75
$ws.Range("I77").Value = 6410
$ws.Range("J77").Value = 100000
$ws.Range("K77").Value = 32050
$ws.Range("L77").Value = 500000
$ws.Range("M77").Value = -27370
$ws.Range("H79").Value = 4477.909
$ws.Range("I79").Value = 4470
$ws.Range("J79").Value = 4499
$ws.Range("K79").Value = 4470
$ws.Range("L79").Value = 4499
$ws.Range("M79").Value = -3378
$ws.Range("N79").Value = -6683
$ws.Range("H80").Value = 1157.6666
$ws.Range("I80").Value = 716
$ws.Range("J80").Value = 1599.3334
$ws.Range("K80").Value = 2148
$ws.Range("L80").Value = 4798.0002
$ws.Range("M80").Value = -1150
$ws.Range("N80").Value = -6794.0002
$ws.Range("H83").Value = 1157.6666
$ws.Range("I83").Value = 716
$ws.Range("J83").Value = 1599.3334
$ws.Range("K83").Value = 6444
$ws.Range("L83").Value = 14394.0006
$ws.Range("M83").Value = -1452
$ws.Range("N83").Value = -24378.0006
$ws.Range("H86").Value = 2423.8572
$ws.Range("I86").Value = 2080.2856
$ws.Range("J86").Value = 2767.4285
$ws.Range("K86").Value = 2080.2856
$ws.Range("L86").Value = 2767.4285
$ws.Range("M86").Value = -957.2856000000002
$ws.Range("N86").Value = -5013.4285
$ws.Range("H88").Value = 5446.1875
$ws.Range("I88").Value = 1126
$ws.Range("J88").Value = 6443.154
$ws.Range("K88").Value = 1126
$ws.Range("L88").Value = 6443.154
$ws.Range("M88").Value = -720
$ws.Range("N88").Value = -7255.154
$ws.Range("H89").Value = 2423.8572
$ws.Range("I89").Value = 2080.2856
$ws.Range("J89").Value = 2767.4285
$ws.Range("K89").Value = 10401.428
$ws.Range("L89").Value = 13837.1425
$ws.Range("M89").Value = -4785.428
$ws.Range("N89").Value = -25069.1425
$ws.Range("H91").Value = 5446.1875
$ws.Range("I91").Value = 1126
$ws.Range("J91").Value = 6443.154
$ws.Range("K91").Value = 1126
$ws.Range("L91").Value = 6443.154
$ws.Range("M91").Value = 278
$ws.Range("N91").Value = -9251.154
$ws.Range("H92").Value = 46215.453
$ws.Range("I92").Value = 50807
$ws.Range("J92").Value = 300
$ws.Range("K92").Value = 50807
$ws.Range("L92").Value = 300
$ws.Range("M92").Value = -49559
$ws.Range("H98").Value = 1830.52
$ws.Range("I98").Value = 1013.45
$ws.Range("J98").Value = 5098.8
$ws.Range("K98").Value = 1013.45
$ws.Range("L98").Value = 5098.8
$ws.Range("M98").Value = 484.55
$ws.Range("H99").Value = 3999.6
$ws.Range("I99").Value = 333
$ws.Range("J99").Value = 9499.5
$ws.Range("K99").Value = 999
$ws.Range("L99").Value = 28498.5
$ws.Range("M99").Value = 499
$ws.Range("N99").Value = -31494.5
$ws.Range("H107").Value = 1487.7273
$ws.Range("I107").Value = 854
$ws.Range("J107").Value = 2596.75
$ws.Range("K107").Value = 854
$ws.Range("L107").Value = 2596.75
$ws.Range("M107").Value = 1066
$ws.Range("N107").Value = -6436.75
$ws.Range("H113").Value = 10080.417
$ws.Range("I113").Value = 19405
$ws.Range("J113").Value = 5418.125
$ws.Range("K113").Value = 19405
$ws.Range("L113").Value = 5418.125
$ws.Range("M113").Value = -16151
$ws.Range("H122").Value = 1830.52
$ws.Range("I122").Value = 1013.45
$ws.Range("J122").Value = 5098.8
$ws.Range("K122").Value = 3040.35
$ws.Range("L122").Value = 15296.4
$ws.Range("M122").Value = -590.3500000000004
$ws.Range("H131").Value = 14657.363
$ws.Range("I131").Value = 4165
$ws.Range("J131").Value = 37141
$ws.Range("K131").Value = 12495
$ws.Range("L131").Value = 111423
$ws.Range("M131").Value = -7455
$ws.Range("N131").Value = -121503
$ws.Range("H132").Value = 8164.2085
$ws.Range("I132").Value = 6348.1
$ws.Range("J132").Value = 17244.75
$ws.Range("K132").Value = 19044.3
$ws.Range("L132").Value = 51734.25
$ws.Range("M132").Value = -16514.3
$ws.Range("H135").Value = 745.5714
$ws.Range("I135").Value = 786.5
$ws.Range("J135").Value = 500
$ws.Range("K135").Value = 7078.5
$ws.Range("L135").Value = 4500
$ws.Range("M135").Value = -4543.5
$ws.Range("N135").Value = -9570
$ws.Range("H137").Value = 1074.697
$ws.Range("I137").Value = 989.5
$ws.Range("J137").Value = 1301.8889
$ws.Range("K137").Value = 2968.5
$ws.Range("L137").Value = 3905.6667
$ws.Range("M137").Value = -418.5
$ws.Range("N137").Value = -9005.6667
$ws.Range("H138").Value = 2864.4
$ws.Range("I138").Value = 2647.3157
$ws.Range("J138").Value = 2965
$ws.Range("K138").Value = 7941.9471
$ws.Range("L138").Value = 8895
$ws.Range("M138").Value = -2801.9471
$ws.Range("N138").Value = -19175
$ws.Range("H141").Value = 2900.0605
$ws.Range("I141").Value = 1777.5385
$ws.Range("J141").Value = 7069.4287
$ws.Range("K141").Value = 5332.6155
$ws.Range("L141").Value = 21208.2861
$ws.Range("M141").Value = -152.6154999999999
$ws.Range("N141").Value = -31568.2861

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2210.0417
$ws.Range("I2").Value = 1559.5555
$ws.Range("J2").Value = 4161.5
$ws.Range("K2").Value = 1559.5555
$ws.Range("L2").Value = 4161.5
$ws.Range("M2").Value = -1446.5555
$ws.Range("H4").Value = 483.33334
$ws.Range("I4").Value = 483.33334
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 483.33334
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -367.33334
$ws.Range("H32").Value = 29305.334
$ws.Range("I32").Value = 41997.4
$ws.Range("J32").Value = 6640.9287
$ws.Range("K32").Value = 41997.4
$ws.Range("L32").Value = 6640.9287
$ws.Range("M32").Value = -41710.4
$ws.Range("H37").Value = 60000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 60000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 60000
$ws.Range("N37").Value = -60546
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H59").Value = 70060
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 70060
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 70060
$ws.Range("N59").Value = -71668
$ws.Range("H61").Value = 3300.36
$ws.Range("I61").Value = 2983.7058
$ws.Range("J61").Value = 3973.25
$ws.Range("K61").Value = 2983.7058
$ws.Range("L61").Value = 3973.25
$ws.Range("M61").Value = -2771.7058
$ws.Range("H74").Value = 3099.5715
$ws.Range("I74").Value = 2884.6956
$ws.Range("J74").Value = 4088
$ws.Range("K74").Value = 2884.6956
$ws.Range("L74").Value = 4088
$ws.Range("M74").Value = -2010.6956
$ws.Range("H77").Value = 3099.5715
$ws.Range("I77").Value = 2884.6956
$ws.Range("J77").Value = 4088
$ws.Range("K77").Value = 14423.478
$ws.Range("L77").Value = 20440
$ws.Range("M77").Value = -10055.478
$ws.Range("H116").Value = 2210.0417
$ws.Range("I116").Value = 1559.5555
$ws.Range("J116").Value = 4161.5
$ws.Range("K116").Value = 1559.5555
$ws.Range("L116").Value = 4161.5
$ws.Range("M116").Value = 734.4445000000001
$ws.Range("H117").Value = 224999.5
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 224999.5
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 224999.5
$ws.Range("N117").Value = -234177.5
$ws.Range("H122").Value = 3197.875
$ws.Range("I122").Value = 2895
$ws.Range("J122").Value = 3298.8333
$ws.Range("K122").Value = 8685
$ws.Range("L122").Value = 9896.499899999999
$ws.Range("M122").Value = -6235
$ws.Range("N122").Value = -14796.4999
$ws.Range("H132").Value = 79063.71000000001
$ws.Range("I132").Value = 205279.8
$ws.Range("J132").Value = 8943.666999999999
$ws.Range("K132").Value = 615839.3999999999
$ws.Range("L132").Value = 26831.001
$ws.Range("M132").Value = -613309.3999999999
$ws.Range("H136").Value = 3300.36
$ws.Range("I136").Value = 2983.7058
$ws.Range("J136").Value = 3973.25
$ws.Range("K136").Value = 8951.117400000001
$ws.Range("L136").Value = 11919.75
$ws.Range("M136").Value = -6401.117400000001

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2210.0417
$ws.Range("I3").Value = 1559.5555
$ws.Range("J3").Value = 4161.5
$ws.Range("K3").Value = 1559.5555
$ws.Range("L3").Value = 4161.5
$ws.Range("M3").Value = -1445.5555
$ws.Range("H36").Value = 41
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 41
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 41
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -1109
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H69").Value = 56666.668
$ws.Range("I69").Value = 55000
$ws.Range("J69").Value = 57500
$ws.Range("K69").Value = 55000
$ws.Range("L69").Value = 57500
$ws.Range("M69").Value = -54189
$ws.Range("N69").Value = -59122
$ws.Range("H72").Value = 56666.668
$ws.Range("I72").Value = 55000
$ws.Range("J72").Value = 57500
$ws.Range("K72").Value = 165000
$ws.Range("L72").Value = 172500
$ws.Range("M72").Value = -160944
$ws.Range("N72").Value = -180612
$ws.Range("H86").Value = 4089.6667
$ws.Range("I86").Value = 2912.1428
$ws.Range("J86").Value = 4678.4287
$ws.Range("K86").Value = 2912.1428
$ws.Range("L86").Value = 4678.4287
$ws.Range("M86").Value = -1789.1428
$ws.Range("N86").Value = -6924.4287
$ws.Range("H89").Value = 4089.6667
$ws.Range("I89").Value = 2912.1428
$ws.Range("J89").Value = 4678.4287
$ws.Range("K89").Value = 14560.714
$ws.Range("L89").Value = 23392.1435
$ws.Range("M89").Value = -8944.714
$ws.Range("N89").Value = -34624.14350000001
$ws.Range("H94").Value = 3264.6667
$ws.Range("I94").Value = 3117.6
$ws.Range("J94").Value = 4000
$ws.Range("K94").Value = 3117.6
$ws.Range("L94").Value = 4000
$ws.Range("M94").Value = -2666.6
$ws.Range("N94").Value = -4902
$ws.Range("H105").Value = 3911.375
$ws.Range("I105").Value = 4140
$ws.Range("J105").Value = 2311
$ws.Range("K105").Value = 4140
$ws.Range("L105").Value = 2311
$ws.Range("M105").Value = -2393
$ws.Range("N105").Value = -5805
$ws.Range("H107").Value = 1831.2
$ws.Range("I107").Value = 1189
$ws.Range("J107").Value = 4400
$ws.Range("K107").Value = 1189
$ws.Range("L107").Value = 4400
$ws.Range("M107").Value = 731
$ws.Range("H134").Value = 2064.6072
$ws.Range("I134").Value = 2046.5385
$ws.Range("J134").Value = 2299.5
$ws.Range("K134").Value = 6139.6155
$ws.Range("L134").Value = 6898.5
$ws.Range("M134").Value = -3604.6155
$ws.Range("N134").Value = -11968.5

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 79999
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 79999
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 79999
$ws.Range("N20").Value = -80471
$ws.Range("H30").Value = 79999
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 79999
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 79999
$ws.Range("N30").Value = -80181
$ws.Range("H31").Value = 2099.56
$ws.Range("I31").Value = 2070.8667
$ws.Range("J31").Value = 2142.6
$ws.Range("K31").Value = 2070.8667
$ws.Range("L31").Value = 2142.6
$ws.Range("M31").Value = -1775.8667
$ws.Range("N31").Value = -2732.6
$ws.Range("H34").Value = 2099.56
$ws.Range("I34").Value = 2070.8667
$ws.Range("J34").Value = 2142.6
$ws.Range("K34").Value = 2070.8667
$ws.Range("L34").Value = 2142.6
$ws.Range("M34").Value = -1868.8667
$ws.Range("N34").Value = -2546.6
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H58").Value = 85599.38
$ws.Range("I58").Value = 256949.5
$ws.Range("J58").Value = 9443.777
$ws.Range("K58").Value = 256949.5
$ws.Range("L58").Value = 9443.777
$ws.Range("M58").Value = -256746.5
$ws.Range("N58").Value = -9849.777
$ws.Range("H62").Value = 3626
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3626
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3626
$ws.Range("N62").Value = -4874
$ws.Range("H65").Value = 3626
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3626
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 18130
$ws.Range("N65").Value = -24370
$ws.Range("H99").Value = 2359.3333
$ws.Range("I99").Value = 2312.4443
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 2312.4443
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -814.4443000000001
$ws.Range("N99").Value = -5496
$ws.Range("H105").Value = 2851.818
$ws.Range("I105").Value = 2056.923
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 2056.923
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -309.9229999999998
$ws.Range("N105").Value = -7494
$ws.Range("H122").Value = 4674
$ws.Range("I122").Value = 2011.5
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 6034.5
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -3584.5
$ws.Range("N122").Value = -34897
$ws.Range("H126").Value = 2359.3333
$ws.Range("I126").Value = 2312.4443
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 6937.3329
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -4467.3329
$ws.Range("N126").Value = -12440
$ws.Range("H128").Value = 79999
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 79999
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 79999
$ws.Range("N128").Value = -89959
$ws.Range("H130").Value = 30000
$ws.Range("I130").Value = 30000
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 30000
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = -24980
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 2664
$ws.Range("I132").Value = 2664
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7992
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5462
$ws.Range("H134").Value = 65897.06
$ws.Range("I134").Value = 70090.2
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 210270.6
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -207735.6
$ws.Range("H136").Value = 85599.38
$ws.Range("I136").Value = 256949.5
$ws.Range("J136").Value = 9443.777
$ws.Range("K136").Value = 770848.5
$ws.Range("L136").Value = 28331.331
$ws.Range("M136").Value = -768298.5
$ws.Range("N136").Value = -33431.331

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 719246.4399999999
$ws.Range("I4").Value = 719246.4399999999
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2157739.32
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -2157627.32
$ws.Range("H13").Value = 14308.143
$ws.Range("I13").Value = 14308.143
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 42924.429
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -42756.429
$ws.Range("H14").Value = 3741
$ws.Range("I14").Value = 3741
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 11223
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -11050
$ws.Range("H17").Value = 55.375
$ws.Range("I17").Value = 47.090908
$ws.Range("J17").Value = 73.59999999999999
$ws.Range("K17").Value = 141.272724
$ws.Range("L17").Value = 220.8
$ws.Range("M17").Value = 27.72727600000002
$ws.Range("N17").Value = -558.8
$ws.Range("H18").Value = 981.44446
$ws.Range("I18").Value = 981.44446
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 2944.33338
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -2775.33338
$ws.Range("H26").Value = 207.66667
$ws.Range("I26").Value = 159.4
$ws.Range("J26").Value = 449
$ws.Range("K26").Value = 478.2
$ws.Range("L26").Value = 1347
$ws.Range("M26").Value = -190.2
$ws.Range("N26").Value = -1923
$ws.Range("H39").Value = 4665
$ws.Range("I39").Value = 592
$ws.Range("J39").Value = 4868.65
$ws.Range("K39").Value = 1776
$ws.Range("L39").Value = 14605.95
$ws.Range("M39").Value = -1482
$ws.Range("N39").Value = -15193.95
$ws.Range("H68").Value = 3166.3333
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 3749.5
$ws.Range("K68").Value = 6000
$ws.Range("L68").Value = 11248.5
$ws.Range("M68").Value = -5189
$ws.Range("N68").Value = -12870.5
$ws.Range("H70").Value = 10538.667
$ws.Range("I70").Value = 7308
$ws.Range("J70").Value = 17000
$ws.Range("K70").Value = 21924
$ws.Range("L70").Value = 51000
$ws.Range("M70").Value = -21609
$ws.Range("N70").Value = -51630
$ws.Range("H71").Value = 3166.3333
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 3749.5
$ws.Range("K71").Value = 18000
$ws.Range("L71").Value = 33745.5
$ws.Range("M71").Value = -13944
$ws.Range("N71").Value = -41857.5
$ws.Range("H73").Value = 10538.667
$ws.Range("I73").Value = 7308
$ws.Range("J73").Value = 17000
$ws.Range("K73").Value = 21924
$ws.Range("L73").Value = 51000
$ws.Range("M73").Value = -20832
$ws.Range("N73").Value = -53184
$ws.Range("H82").Value = 5698.4287
$ws.Range("I82").Value = 4472.25
$ws.Range("J82").Value = 7333.3335
$ws.Range("K82").Value = 13416.75
$ws.Range("L82").Value = 22000.0005
$ws.Range("M82").Value = -13010.75
$ws.Range("N82").Value = -22812.0005
$ws.Range("H85").Value = 5698.4287
$ws.Range("I85").Value = 4472.25
$ws.Range("J85").Value = 7333.3335
$ws.Range("K85").Value = 13416.75
$ws.Range("L85").Value = 22000.0005
$ws.Range("M85").Value = -12012.75
$ws.Range("N85").Value = -24808.0005
$ws.Range("H86").Value = 787.375
$ws.Range("I86").Value = 785.5714
$ws.Range("J86").Value = 800
$ws.Range("K86").Value = 2356.7142
$ws.Range("L86").Value = 2400
$ws.Range("M86").Value = -1170.7142
$ws.Range("N86").Value = -4772
$ws.Range("H89").Value = 787.375
$ws.Range("I89").Value = 785.5714
$ws.Range("J89").Value = 800
$ws.Range("K89").Value = 7070.1426
$ws.Range("L89").Value = 7200
$ws.Range("M89").Value = -1142.1426
$ws.Range("N89").Value = -19056
$ws.Range("H92").Value = 507.26666
$ws.Range("I92").Value = 277.84616
$ws.Range("J92").Value = 1998.5
$ws.Range("K92").Value = 833.5384799999999
$ws.Range("L92").Value = 5995.5
$ws.Range("M92").Value = 414.4615200000001
$ws.Range("N92").Value = -8491.5
$ws.Range("H131").Value = 20630.188
$ws.Range("I131").Value = 20809.2
$ws.Range("J131").Value = 20548.818
$ws.Range("K131").Value = 62427.60000000001
$ws.Range("L131").Value = 61646.454
$ws.Range("M131").Value = -57387.60000000001
$ws.Range("N131").Value = -71726.454

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4447.0527
$ws.Range("I70").Value = 4120.3
$ws.Range("J70").Value = 4810.1113
$ws.Range("K70").Value = 4120.3
$ws.Range("L70").Value = 4810.1113
$ws.Range("M70").Value = -3850.3
$ws.Range("H73").Value = 4447.0527
$ws.Range("I73").Value = 4120.3
$ws.Range("J73").Value = 4810.1113
$ws.Range("K73").Value = 4120.3
$ws.Range("L73").Value = 4810.1113
$ws.Range("M73").Value = -3184.3
$ws.Range("H97").Value = 952.1739
$ws.Range("I97").Value = 761.5
$ws.Range("J97").Value = 1638.6
$ws.Range("K97").Value = 761.5
$ws.Range("L97").Value = 1638.6
$ws.Range("M97").Value = -265.5
$ws.Range("H122").Value = 4118
$ws.Range("I122").Value = 3210
$ws.Range("J122").Value = 5285.4287
$ws.Range("K122").Value = 9630
$ws.Range("L122").Value = 15856.2861
$ws.Range("M122").Value = -7180
$ws.Range("N122").Value = -20756.2861
$ws.Range("H123").Value = 57329.668
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 57329.668
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 57329.668
$ws.Range("N123").Value = -62229.668
$ws.Range("H130").Value = 89999.5
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 89999.5
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 89999.5
$ws.Range("N130").Value = -100039.5
$ws.Range("H132").Value = 50261.363
$ws.Range("I132").Value = 67346.94
$ws.Range("J132").Value = 4699.8335
$ws.Range("K132").Value = 202040.82
$ws.Range("L132").Value = 14099.5005
$ws.Range("M132").Value = -199510.82
$ws.Range("N132").Value = -19159.5005

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I22").Value = 93276.414
$ws.Range("J22").Value = 1423.4445
$ws.Range("K22").Value = 93276.414
$ws.Range("L22").Value = 1423.4445
$ws.Range("M22").Value = -92981.414
$ws.Range("N22").Value = -2013.4445
$ws.Range("I27").Value = 93276.414
$ws.Range("J27").Value = 1423.4445
$ws.Range("K27").Value = 93276.414
$ws.Range("L27").Value = 1423.4445
$ws.Range("M27").Value = -93169.414
$ws.Range("N27").Value = -1637.4445
$ws.Range("H40").Value = 5037.375
$ws.Range("I40").Value = 4792.7856
$ws.Range("J40").Value = 6749.5
$ws.Range("K40").Value = 4792.7856
$ws.Range("L40").Value = 6749.5
$ws.Range("M40").Value = -4656.7856
$ws.Range("N40").Value = -7021.5
$ws.Range("H46").Value = 6703.28
$ws.Range("I46").Value = 37799.668
$ws.Range("J46").Value = 2462.8635
$ws.Range("K46").Value = 37799.668
$ws.Range("L46").Value = 2462.8635
$ws.Range("M46").Value = -37611.668
$ws.Range("N46").Value = -2838.8635
$ws.Range("H68").Value = 3578.7646
$ws.Range("I68").Value = 2217.8
$ws.Range("J68").Value = 5523
$ws.Range("K68").Value = 2217.8
$ws.Range("L68").Value = 5523
$ws.Range("M68").Value = -1468.8
$ws.Range("N68").Value = -7021
$ws.Range("H71").Value = 3578.7646
$ws.Range("I71").Value = 2217.8
$ws.Range("J71").Value = 5523
$ws.Range("K71").Value = 11089
$ws.Range("L71").Value = 27615
$ws.Range("M71").Value = -7345
$ws.Range("N71").Value = -35103
$ws.Range("H122").Value = 4812.5
$ws.Range("I122").Value = 3500
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 10500
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -8050
$ws.Range("H134").Value = 94999
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 94999
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 94999
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -105139
$ws.Range("H136").Value = 2328.25
$ws.Range("I136").Value = 2031.6111
$ws.Range("J136").Value = 4998
$ws.Range("K136").Value = 6094.8333
$ws.Range("L136").Value = 14994
$ws.Range("M136").Value = -3544.8333

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2298.6667
$ws.Range("I96").Value = 1967.3334
$ws.Range("J96").Value = 3292.6667
$ws.Range("K96").Value = 1967.3334
$ws.Range("L96").Value = 3292.6667
$ws.Range("M96").Value = -594.3334
$ws.Range("N96").Value = -6038.6667
$ws.Range("H100").Value = 956.36365
$ws.Range("I100").Value = 874.4
$ws.Range("J100").Value = 1132
$ws.Range("K100").Value = 1748.8
$ws.Range("L100").Value = 2264
$ws.Range("M100").Value = -1207.8
$ws.Range("H122").Value = 9172.6
$ws.Range("I122").Value = 8964.5
$ws.Range("J122").Value = 10005
$ws.Range("K122").Value = 26893.5
$ws.Range("L122").Value = 30015
$ws.Range("M122").Value = -24443.5
$ws.Range("H126").Value = 80286.46000000001
$ws.Range("I126").Value = 93978.63
$ws.Range("J126").Value = 4979.5
$ws.Range("K126").Value = 281935.89
$ws.Range("L126").Value = 14938.5
$ws.Range("M126").Value = -279465.89
$ws.Range("H132").Value = 102717.7
$ws.Range("I132").Value = 113575.22
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 340725.66
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -338195.66
$ws.Range("H136").Value = 6051.609
$ws.Range("I136").Value = 6921.1333
$ws.Range("J136").Value = 4421.25
$ws.Range("K136").Value = 20763.3999
$ws.Range("L136").Value = 13263.75
$ws.Range("M136").Value = -18213.3999
$ws.Range("H141").Value = 85684.39999999999
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 85684.39999999999
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 85684.39999999999
$ws.Range("N141").Value = -96044.39999999999
